# PagoAnticipadoValor.xlsx - "update entregable 1, 2"
#
# The sheet had 8 columns (A:H): usuario, contraseña, arreglo prestamo,
# cuenta debito, usuarioAp, Estado, Transaccion, Fecha - with sample
# values in row 2. The edit removes the "cuenta debito" (D) and
# "usuarioAp" (E) columns entirely (shifting Estado/Transaccion/Fecha
# left into D:F), changes the "usuario" sample value in A2, and clears
# the now-unused trailing sample values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove columns D ("cuenta debito") and E ("usuarioAp") completely;
# F/G/H (Estado/Transaccion/Fecha headers + PASSED/transaction id/date
# samples) shift left to become D/E/F.
$ws.Range("D1:E1").EntireColumn.Delete()

# The sample "usuario" value in A2 changes from SCISNEROSC1 to ebenito.
$ws.Cells.Item(2, 1).Value2 = "ebenito"

# The old F2/G2/H2 sample values (PASSED / PI000000139885 20 /
# 16 abr. 2023, 20:14:51) are no longer present after the shift.
$ws.Cells.Item(2, 4).Value2 = ""
$ws.Cells.Item(2, 5).Value2 = ""
$ws.Cells.Item(2, 6).Value2 = ""

# Match the saved selection/active cell.
$ws.Range("C8").Select()
